# "Drop in RMI script results for 3.0"
#
# This restores the original RMI-sourced hydrogen "percent excess capacity"
# figures, undoing a Texas-specific customization:
#   1. Remove the "Data Texas" worksheet (and the long commentary notes that
#      lived only on that sheet - deleting the sheet drops those shared
#      strings automatically).
#   2. Reset the input assumption on the HPPECbP sheet (cell B2) from 10%
#      back to the RMI study's 25% excess capacity figure. All the other
#      cells on that sheet are formulas referencing B2 (directly or via the
#      shared "=$B$2" formula), so they recalculate automatically.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Drop the "Data Texas" sheet entirely.
$dataTexas = $wb.Worksheets.Item("Data Texas")
$dataTexas.Delete()

# 2. Update the excess-capacity assumption back to 25%.
$hppecbp = $wb.Worksheets.Item("HPPECbP")
$hppecbp.Range("B2").Value = 0.25

$excel.DisplayAlerts = $true
